# "start before cleaning the data"
# Adds a new raw/unclean data block (rows 27-49) to the "kliceni" sheet:
# a small header row followed by per-sample seed-count rows, all dated
# 42016 (2015-01-12). Also refreshes the column widths on both sheets and
# moves the active selection to reflect where editing left off.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("kliceni")
$ws2 = $wb.Worksheets.Item("proklicovani")

# --- New header row 27 -------------------------------------------------
# Values are entered in this order so the shared-string table gets the
# same new-entry ordering (seed_n, N, dish, locality, date) as the target.
$ws1.Range("D27").Value = "seed_n"
$ws1.Range("E27").Value = "N"
$ws1.Range("C27").Value = "dish"
$ws1.Range("F27").Value = "locality"
$ws1.Range("B27").Value = "date"

# --- New data rows 28-49 ------------------------------------------------
$labels = @("os1","os2","os3","ka1","ka2","ka3","no1","st1","st2","st3", `
            "ry1","ry2","ry3","kr1","kr2","kr3","mr1","mr2","mr3","la1", `
            "la2","la3")
$dvals  = @(1,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 28 + $i
    $ws1.Range("B$row").Value = 42016
    $ws1.Range("C$row").Value = $labels[$i]
    $ws1.Range("D$row").Value = $dvals[$i]
}

# Give the new date column (B28:B49) the same date formatting as the
# existing header date cells (B1:Q1) by copying their format only.
$ws1.Range("B1").Copy()
$ws1.Range("B28:B49").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths on both sheets --------------------------------------
$ws1.Columns("B:Q").ColumnWidth = 10.140625
$ws2.Columns("B:Q").ColumnWidth = 10.140625

# --- Selections reflecting where work stopped ---------------------------
$ws1.Range("F39").Select()
$ws2.Range("D20").Select()
